# Weekly fruit/vegetable price update: a new price record is inserted at
# the top of the Espinaca price-history block (row 263), pushing all the
# existing historical rows (old 263..296) down by one row (new 264..297).
#
# $ws is the (only) worksheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 263 - this shifts rows 263:296 down to
# 264:297 and grows the sheet's used range accordingly (R296 -> R297).
$ws.Rows.Item(263).Insert()

# Populate the newly-inserted row 263 with the new weekly record. All the
# "descriptive" columns (market, region, product, quality, unit, origin,
# classification, etc.) repeat the same values used throughout this
# sub-block; only the date and the price/volume figures are new.
$ws.Cells.Item(263, 1).Value = 8
$ws.Cells.Item(263, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(263, 3).Value = "Coquimbo"
$ws.Cells.Item(263, 4).Value = 44776
$ws.Cells.Item(263, 5).Value = 4
$ws.Cells.Item(263, 6).Value = 100112012
$ws.Cells.Item(263, 7).Value = "Espinaca"
$ws.Cells.Item(263, 8).Value = "Sin especificar"
$ws.Cells.Item(263, 9).Value = "Primera"
$ws.Cells.Item(263, 10).Value = 2800
$ws.Cells.Item(263, 11).Value = 550
$ws.Cells.Item(263, 12).Value = 600
$ws.Cells.Item(263, 13).Value = 575
$ws.Cells.Item(263, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(263, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(263, 16).Value = 1150
$ws.Cells.Item(263, 17).Value = 0.5
$ws.Cells.Item(263, 18).Value = "Hortaliza"

# Match the date number format used by the other "Fecha" cells in column D.
$ws.Cells.Item(263, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
